# Auto-generated edit script: updates crypto price/volume/coin data for Feb 13 2023 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Prefix with an apostrophe so Excel stores the literal text instead of
    # auto-converting number-/percent-looking strings, then reset the style
    # back to Normal so the quote-prefix indicator does not leave a style diff.
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# Row 2
Set-TextValue "D2" '299.45'
Set-TextValue "E2" '-2.74%'
# Row 3
Set-TextValue "D3" '40.39'
Set-TextValue "E3" '-1.61%'
# Row 4
Set-TextValue "D4" '5.154'
Set-TextValue "E4" '-1.11%'
# Row 5
Set-TextValue "D5" '0.07515'
Set-TextValue "E5" '-1.99%'
# Row 6
Set-TextValue "D6" '4.343'
Set-TextValue "E6" '0.81%'
# Row 7
Set-TextValue "D7" '1.627'
Set-TextValue "E7" '-0.36%'
# Row 8
Set-TextValue "D8" '0.9393'
Set-TextValue "E8" '2.63%'
# Row 9
Set-TextValue "B9" 'BTSEToken'
Set-TextValue "C9" 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue "D9" '2.419'
Set-TextValue "E9" '-0.41%'
# Row 10
Set-TextValue "B10" 'LiechtensteinCryptoassetsExchange'
Set-TextValue "C10" 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue "D10" '0.1197'
Set-TextValue "E10" '-1.14%'
# Row 11
Set-TextValue "B11" 'WazirX'
Set-TextValue "C11" 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue "D11" '0.1796'
Set-TextValue "E11" '-1.27%'
# Row 12
Set-TextValue "B12" 'MandalaExchangeToken'
Set-TextValue "C12" 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue "D12" '0.08957'
Set-TextValue "E12" '-2.07%'
# Row 13
Set-TextValue "B13" 'BitrueCoin'
Set-TextValue "C13" 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue "D13" '0.04145'
Set-TextValue "E13" '-0.32%'
# Row 14
Set-TextValue "B14" 'BitMartToken'
Set-TextValue "C14" 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue "D14" '0.1055'
Set-TextValue "E14" '0.42%'
# Row 15
Set-TextValue "B15" 'BitForexToken'
Set-TextValue "C15" 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue "D15" '0.001263'
Set-TextValue "E15" '0.20%'
# Row 16
Set-TextValue "B16" 'TigerCash'
Set-TextValue "C16" 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue "D16" '0.005862'
Set-TextValue "E16" '1.61%'
# Row 17
Set-TextValue "B17" 'LEO'
Set-TextValue "C17" 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue "D17" '3.347'
Set-TextValue "E17" '0.18%'
# Row 18
Set-TextValue "D18" '0.3352'
Set-TextValue "E18" '0.49%'
# Row 19
Set-TextValue "D19" '7.636'
Set-TextValue "E19" '3.27%'
# Row 20
Set-TextValue "D20" '0.1357'
Set-TextValue "E20" '-1.91%'
# Row 22
Set-TextValue "D22" '0.03838'
Set-TextValue "E22" '-4.46%'
# Row 23
Set-TextValue "D23" '0.001289'
Set-TextValue "E23" '2.44%'
# Row 24
Set-TextValue "D24" '0.003964'
Set-TextValue "E24" '-9.51%'
# Row 25
Set-TextValue "D25" '0.0001306'
Set-TextValue "E25" '0.45%'
# Row 26
Set-TextValue "D26" '0.0003739'
Set-TextValue "E26" '-95.02%'
# Row 38
Set-TextValue "D38" '0.02361'
Set-TextValue "E38" '-5.23%'
# Row 39
Set-TextValue "D39" '0.05080'
Set-TextValue "E39" '-4.67%'
# Row 40
Set-TextValue "D40" '0.007759'
Set-TextValue "E40" '-0.96%'
# Row 41
Set-TextValue "D41" '0.1299'
Set-TextValue "E41" '-1.01%'
# Row 42
Set-TextValue "D42" '0.007602'
Set-TextValue "E42" '16.84%'
# Row 43
Set-TextValue "D43" '0.003569'
Set-TextValue "E43" '86.87%'
# Row 44
Set-TextValue "D44" '0.007401'
Set-TextValue "E44" '-10.28%'
# Row 45
Set-TextValue "D45" '0.3246'
Set-TextValue "E45" '-2.77%'
# Row 46
Set-TextValue "D46" '0.00006824'
Set-TextValue "E46" '1.77%'
# Row 47
Set-TextValue "E47" '0.30%'
# Row 48
Set-TextValue "D48" '0.2463'
Set-TextValue "E48" '-20.40%'
# Row 49
Set-TextValue "D49" '0.004219'
Set-TextValue "E49" '35.90%'
# Row 50
Set-TextValue "D50" '0.00002109'
Set-TextValue "E50" '0.30%'
# Row 51
Set-TextValue "D51" '0.0002009'
Set-TextValue "E51" '0.30%'
